# The "Acc_Upfront" sheet actually contains three separate blocks of data
# (rows 2-3, rows 5-6 and rows 8-9, each preceded by the same header row 1
# and separated by a blank spacer row). This edit splits those three
# blocks out into three individual worksheets - Acc_Upfront1 (block
# 140/141), Acc_Upfront2 (block 142/143) and Acc_Upfront3 (block
# 149/150, which simply stays in the original sheet/tab, renamed) - all
# placed immediately after the original "Acc_Upfront" position, in that
# tab order.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Acc_Upfront")

# --- New sheet: Acc_Upfront1 -> block 1 (IDs 140/141, originally rows 2:3)
$ws1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $src)
$ws1.Name = "Acc_Upfront1"
$src.Range("A1:I3").Copy($ws1.Range("A1"))
$ws1.Range("D2").Select()

# --- New sheet: Acc_Upfront2 -> block 2 (IDs 142/143, originally rows 5:6)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Acc_Upfront2"
$src.Range("A1:G1").Copy($ws2.Range("A1"))
$src.Range("H1:I1").Copy($ws2.Range("H1"))
$src.Range("A5:I6").Copy($ws2.Range("A2"))
$ws2.Range("G3").Select()

# --- Original sheet keeps block 3 (IDs 149/150, originally rows 8:9) and
# is renamed Acc_Upfront3; everything else (old blocks 1 & 2, plus the
# blank spacer rows) is removed so the remaining data shifts up to
# rows 2:3.
$src.Rows("2:7").Delete()
$src.Name = "Acc_Upfront3"
$src.Activate()
$src.Range("D2").Select()
